$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 - this shifts the existing data rows
# (Vanlig padda / Skogsödla / Björktrast) down from rows 2-4 to rows 3-5,
# preserving their content untouched.
$ws.Rows("2:2").Insert()

# Populate the newly inserted row 2 with the new observation record.
$ws.Range("A2").Value = 89087437
$ws.Range("B2").Value = 88806
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 5685
$ws.Range("F2").Value = "Gullgröppa"
$ws.Range("G2").Value = "Pseudomerulius aureus"
$ws.Range("H2").Value = "(Fr.) Jülich"

# Empty (present-but-blank) cells, matching the sparse layout of the new row.
# A plain Value="" assignment removes the cell entirely, so nudge the engine
# into keeping a blank placeholder cell by touching its number format instead.
$ws.Range("I2").NumberFormat = "General"
$ws.Range("J2").NumberFormat = "General"
$ws.Range("K2").NumberFormat = "General"
$ws.Range("N2").NumberFormat = "General"

$ws.Range("P2").Value = "Kåddis, Vb"
$ws.Range("Q2").Value = 750788.3443989656
$ws.Range("R2").Value = 7089862.722689836
$ws.Range("S2").Value = 25
$ws.Range("T2").Value = "Västerbotten"
$ws.Range("U2").Value = "Umeå"
$ws.Range("V2").Value = "Västerbotten"
$ws.Range("W2").Value = "Umeå socken"

# Force these date-like strings to be stored as plain text (not parsed as dates).
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2020-11-15"
$ws.Range("Z2").Value = "00:00"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2020-11-15"
$ws.Range("AB2").Value = "00:00"

$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AF2").NumberFormat = "General"
$ws.Range("AG2").Value = $false
$ws.Range("AT2").NumberFormat = "General"

$ws.Range("AW2").Value = "Gunhild Eriksson Nyberg"
$ws.Range("AX2").Value = "Gunhild Eriksson Nyberg"
$ws.Range("AY2").NumberFormat = "General"
